$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("a")

# --- Update the cumulative-frequency table a!I1:Z7 ---
# Columns I..U already held values; V..Z previously held #N/A errors
# and are now populated with real numbers as more data points were added.
$data = New-Object 'object[,]' 7,18
$data[0,0] = 169
$data[0,1] = 329.5
$data[0,2] = 491.5
$data[0,3] = 631
$data[0,4] = 776.5
$data[0,5] = 914
$data[0,6] = 1054
$data[0,7] = 1169
$data[0,8] = 1283.5
$data[0,9] = 1402
$data[0,10] = 1516.5
$data[0,11] = 1628
$data[0,12] = 1721
$data[0,13] = 1811
$data[0,14] = 1897.5
$data[0,15] = 1981.5
$data[0,16] = 2059
$data[0,17] = 2140.5

$data[1,0] = 169
$data[1,1] = 329.5
$data[1,2] = 491.5
$data[1,3] = 631
$data[1,4] = 776.5
$data[1,5] = 914
$data[1,6] = 1054
$data[1,7] = 1169
$data[1,8] = 1283.5
$data[1,9] = 1402
$data[1,10] = 1516.5
$data[1,11] = 1628
$data[1,12] = 1721
$data[1,13] = 1811
$data[1,14] = 1897.5
$data[1,15] = 1981.5
$data[1,16] = 2059
$data[1,17] = 2140.5

$data[2,0] = 140
$data[2,1] = 267
$data[2,2] = 381
$data[2,3] = 507.5
$data[2,4] = 628.5
$data[2,5] = 735.5
$data[2,6] = 842
$data[2,7] = 956
$data[2,8] = 1064.5
$data[2,9] = 1152.5
$data[2,10] = 1231.5
$data[2,11] = 1322
$data[2,12] = 1409
$data[2,13] = 1487
$data[2,14] = 1560
$data[2,15] = 1636
$data[2,16] = 1707
$data[2,17] = 1768.5

$data[3,0] = 140
$data[3,1] = 267
$data[3,2] = 381
$data[3,3] = 507.5
$data[3,4] = 628.5
$data[3,5] = 735.5
$data[3,6] = 842
$data[3,7] = 956
$data[3,8] = 1064.5
$data[3,9] = 1152.5
$data[3,10] = 1231.5
$data[3,11] = 1322
$data[3,12] = 1409
$data[3,13] = 1487
$data[3,14] = 1560
$data[3,15] = 1636
$data[3,16] = 1707
$data[3,17] = 1768.5

$data[4,0] = 91.5
$data[4,1] = 192
$data[4,2] = 280.5
$data[4,3] = 366
$data[4,4] = 449.5
$data[4,5] = 536.5
$data[4,6] = 601.5
$data[4,7] = 678.5
$data[4,8] = 755.5
$data[4,9] = 814
$data[4,10] = 876.5
$data[4,11] = 928.5
$data[4,12] = 986.5
$data[4,13] = 1042
$data[4,14] = 1099
$data[4,15] = 1141
$data[4,16] = 1182.5
$data[4,17] = 1224.5

$data[5,0] = 91.5
$data[5,1] = 192
$data[5,2] = 280.5
$data[5,3] = 366
$data[5,4] = 449.5
$data[5,5] = 536.5
$data[5,6] = 601.5
$data[5,7] = 678.5
$data[5,8] = 755.5
$data[5,9] = 814
$data[5,10] = 876.5
$data[5,11] = 928.5
$data[5,12] = 986.5
$data[5,13] = 1042
$data[5,14] = 1099
$data[5,15] = 1141
$data[5,16] = 1182.5
$data[5,17] = 1224.5

$data[6,0] = 437
$data[6,1] = 852
$data[6,2] = 1267
$data[6,3] = 1662
$data[6,4] = 2015
$data[6,5] = 2361
$data[6,6] = 2704
$data[6,7] = 3016
$data[6,8] = 3299
$data[6,9] = 3612
$data[6,10] = 3905
$data[6,11] = 4164
$data[6,12] = 4419
$data[6,13] = 4668
$data[6,14] = 4898
$data[6,15] = 5124
$data[6,16] = 5343
$data[6,17] = 5542

$ws.Range("I1:Z7").Value = $data

# --- Update the cached summary values in a!C10:D15 (plain numbers, not formulas) ---
$ws.Range("C10").Value = 8281
$ws.Range("D10").Value = 68
$ws.Range("C11").Value = 5425.983763200001
$ws.Range("D11").Value = 67
$ws.Range("C14").Value = 10496.481944256
$ws.Range("D14").Value = 88
$ws.Range("C15").Value = 9973.5
$ws.Range("D15").Value = 79

# --- B10:B15 are formulas (+I1, +I2, ...) and a!C12/D12, a!C13/D13 stay at 0; ---
# --- Worksheet "Hoja2" formulas (D3:X9, D10:M10, E14:H19) recalc automatically. ---
$excel.CalculateFullRebuild()
